# "cryptos.xlsx" is refreshed with the latest market snapshot: the Price
# (column D) and Volume(1h) (column E) columns are updated for most rows, and
# rows 47-49 are re-ranked (Coin/Link/Price/Volume shuffled between them).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Maps each changed cell to its new text value.
$updates = [ordered]@{
    'D2' = '39.521.47'
    'E2' = '  +1.93%  '
    'D3' = '2.168.30'
    'E3' = '  +3.11%  '
    'E4' = '  -0.02%  '
    'D5' = '229.03'
    'E5' = '  +0.54%  '
    'D6' = '0.634'
    'E6' = '  +2.84%  '
    'D7' = '63.71'
    'E7' = '  +2.24%  '
    'E8' = '  +0.03%  '
    'D9' = '0.396'
    'E9' = '  +1.49%  '
    'E10' = '  +1.53%  '
    'E11' = '  +0.12%  '
    'E12' = '  +2.20%  '
    'D13' = '2.490.17'
    'E13' = '  +3.10%  '
    'D14' = '22.06'
    'E14' = '  +0.11%  '
    'D15' = '0.814'
    'E15' = '  +0.71%  '
    'D16' = '5.53'
    'E16' = '  -0.05%  '
    'D17' = '2.169.97'
    'E17' = '  +3.05%  '
    'D18' = '39.493.04'
    'E18' = '  +1.89%  '
    'D19' = '6.22'
    'E19' = '  +1.85%  '
    'D20' = '71.94'
    'E20' = '  +0.06%  '
    'E21' = '  +1.05%  '
    'D22' = '229.60'
    'E22' = '  +0.83%  '
    'E23' = '  +0.06%  '
    'E24' = '  +1.84%  '
    'D25' = '2.30'
    'E25' = '  -1.41%  '
    'D26' = '9.74'
    'E26' = '  +1.04%  '
    'D27' = '172.28'
    'E27' = '  -0.03%  '
    'D28' = '0.139'
    'E28' = '  +0.14%  '
    'D29' = '19.91'
    'E29' = '  +2.90%  '
    'E30' = '  +0.20%  '
    'D31' = '2.64'
    'E31' = '  +4.88%  '
    'D32' = '0.124'
    'E32' = '  +2.43%  '
    'E33' = '  +1.70%  '
    'D34' = '4.74'
    'E34' = '  -0.62%  '
    'E35' = '  +0.69%  '
    'E36' = '  +0.01%  '
    'E37' = '  +1.17%  '
    'D38' = '3.62'
    'E38' = '  +1.20%  '
    'E39' = '  +0.20%  '
    'D40' = '102.97'
    'E40' = '  +0.23%  '
    'D41' = '17.96'
    'E41' = '  -1.00%  '
    'D42' = '0.0229'
    'E42' = '  +0.06%  '
    'D43' = '1.525.45'
    'E43' = '  -0.04%  '
    'E44' = '  +0.07%  '
    'E45' = '  +5.46%  '
    'D46' = '4.30'
    'E46' = '  +3.87%  '
    'B47' = 'FraxShare'
    'C47' = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
    'D47' = '7.90'
    'E47' = '  +1.35%  '
    'B48' = 'HuobiToken'
    'C48' = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
    'D48' = '2.82'
    'E48' = '  +0.46%  '
    'B49' = 'Cronos'
    'C49' = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
    'D49' = '0.0924'
    'E49' = '  +1.43%  '
    'D50' = '2.372.25'
    'E50' = '  +3.11%  '
    'E51' = '  -0.66%  '
}

# Cells whose new value looks like a plain number (e.g. "229.03", "0.634").
# Excel would otherwise auto-convert these to numeric cells, losing their
# exact textual form, so for these we briefly force the Text number format
# while writing the value, then clear the formatting again so the cell keeps
# its original (unstyled) look, just like the other text cells.
$forceTextRefs = @(
    'D5'
    'D6'
    'D7'
    'D9'
    'D14'
    'D15'
    'D16'
    'D19'
    'D20'
    'D22'
    'D25'
    'D26'
    'D27'
    'D28'
    'D29'
    'D31'
    'D32'
    'D34'
    'D38'
    'D40'
    'D41'
    'D42'
    'D46'
    'D47'
    'D48'
    'D49'
)

foreach ($cellRef in $updates.Keys) {
    $newValue = $updates[$cellRef]
    $cell = $ws.Range($cellRef)

    if ($forceTextRefs -contains $cellRef) {
        $cell.NumberFormat = "@"
        $cell.Value = $newValue
        $cell.ClearFormats()
    } else {
        $cell.Value = $newValue
    }
}
